$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell C1 - reuse the same formatting as A1/B1 (bold, bordered, centered)
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C1").Value = "Coord: normal vector scan"

# Updated angle values (column B) and new coordinate strings (column C)
$values = @(
    @(0.2243166572380701, "[0.         0.32126117 0.94699063]"),
    @(1.466880383355958, "[-0.43682036  0.52375887  0.7313444 ]"),
    @(0.133319600267756, "[-0.00179395  0.01752325  0.99984485]"),
    @(0.8344448762653219, "[-0.00124692  0.28056321 -0.95983474]"),
    @(2.169621105349656, "[0.71968397 0.28546058 0.63290382]"),
    @(1.166217900643735, "[-0.73138087 -0.26632338  0.62781676]"),
    @(0.5702410547430989, "[0.         0.31553787 0.94891298]"),
    @(0.9908661927219888, "[ 0.         -0.30856319  0.95120385]"),
    @(2.36264967122835, "[-0.72389606  0.25239499  0.64208353]"),
    @(2.134941892813619, "[ 0.72543151 -0.27526146  0.63085676]"),
    @(0.8428997607208314, "[ 0.00125139 -0.28042141 -0.95987617]"),
    @(2.40975649655719, "[-0.71713305 -0.28467602  0.63614444]"),
    @(3.935958454348492, "[0.70126719 0.25677996 0.66504765]"),
    @(1.228143127260918, "[ 0.         -0.30462137  0.95247353]")
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i][0]
    $ws.Cells.Item($row, 3).Value = $values[$i][1]
}
